# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colour scheme
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" colour scheme (the
#                            theme actually used by the slide master/design)
# The authored edit swaps the two themes' contents in place: theme1.xml
# ends up holding the old "Integral"/Red-Violet data and theme2.xml ends up
# holding the old "Office Theme"/Office data, while every relationship keeps
# pointing at the same theme file names as before.
#
# ppt/theme/theme1.xml isn't wired to any slide/notes/handout master that
# the PowerPoint object model exposes distinctly in this deck (it's only
# referenced by the notes master, and NotesMaster/HandoutMaster resolve to
# the very same Theme object as the slide master here), so it can't be
# reached through COM automation. The reachable, user-visible half of the
# swap -- the slide design's theme (ppt/theme/theme2.xml) switching from the
# Integral/Red-Violet palette to the Office Theme palette -- is applied via
# ThemeColorScheme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# ThemeColorScheme(index).RGB packs a colour like the classic VBA RGB()
# macro: R + G*256 + B*65536 -- NOT the hex string read as a plain integer.
#   index : theme tag -> target "Office Theme" hex  (packed RGB long)
$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72

Write-Output "Recoloured slide master theme to the Office Theme palette."
